$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the displayed sale price text cell (was "45.0000", now "60.0000") ---
# This cell is numerically-formatted (0.00) but stores a literal text string,
# so we temporarily force a text number format to avoid Excel re-casting the
# assigned string into a numeric value, then restore the original format.
$p7 = $ws.Range("P7")
$p7Format = $p7.NumberFormat
$p7.NumberFormat = "@"
$p7.Value = "60.0000"
$p7.NumberFormat = $p7Format

# --- Update the "order count" text cell (was "3:0", now "4:0") ---
$ws.Range("Q7").Value = "4:0"

# --- Update the underlying numeric sale price (was 45, now 60) ---
$ws.Range("P8").Value = 60

# --- Update the generated timestamp text (was 10:17 AM, now 10:32 AM) ---
$ws.Range("A9").Value = "Sunday, 14 September, 2025 10:32 AM"
